$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple value replacements (rows 1-4) ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text = "220"

# --- Row 6: 0.00056 -> 0.00249 ---
$t.Rows.Item(6).Cells.Item(1).Range.Text = "0.00249"

# --- Delete rows 7, 8, 9 (0.00014, 0.00004, 0.00017) ---
# delete from highest index to lowest so indices of earlier rows stay valid
$t.Rows.Item(9).Delete()
$t.Rows.Item(8).Delete()
$t.Rows.Item(7).Delete()

# After the deletions, the table collapsed:
#  row 7 -> 0.00019 (unchanged)
#  row 8 -> 0.00023 (needs -> 0.00012)
#  row 9 -> 0.01539 (needs -> 0.00029)
$t.Rows.Item(8).Cells.Item(1).Range.Text = "0.00012"
$t.Rows.Item(9).Cells.Item(1).Range.Text = "0.00029"

# --- Insert 3 new rows after row 9 (before former row 10, "100.0") ---
$newRow1 = $t.Rows.Add($t.Rows.Item(10))
$newRow1.Cells.Item(1).Range.Text = "0.00036"

$newRow2 = $t.Rows.Add($t.Rows.Item(11))
$newRow2.Cells.Item(1).Range.Text = "0.00044"

$newRow3 = $t.Rows.Add($t.Rows.Item(12))
$newRow3.Cells.Item(1).Range.Text = "0.04743"

# --- Rows 44, 45, 46 (tab-separated) -> single condensed value ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.93"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.05"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "65"
